# Weekly update: insert 3 new price rows for "Vega Monumental Concepción - Lechuga"
# right before the existing row that corresponds to date 44545 (old row 1110),
# shifting all subsequent rows down by 3 (old 1110..1185 -> new 1113..1188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 1110, 1111, 1112 (pushes old 1110.. down to 1113..)
$ws.Range("A1110:A1112").EntireRow.Insert()

# --- New row 1110: Conconina(o) ---
$ws.Range("A1110").Value = 11
$ws.Range("B1110").Value = "Vega Monumental Concepción"
$ws.Range("C1110").Value = "Bíobío"
$ws.Range("D1110").Value = 45106
$ws.Range("E1110").Value = 8
$ws.Range("F1110").Value = 100112033
$ws.Range("G1110").Value = "Lechuga"
$ws.Range("H1110").Value = "Conconina(o)"
$ws.Range("I1110").Value = "Primera"
$ws.Range("J1110").Value = 100
$ws.Range("K1110").Value = 6000
$ws.Range("L1110").Value = 6500
$ws.Range("M1110").Value = 6250
$ws.Range("N1110").Value = "$/caja 10 unidades"
$ws.Range("O1110").Value = "Región Metropolitana"
$ws.Range("P1110").Value = 625
$ws.Range("Q1110").Value = 10
$ws.Range("R1110").Value = "Hortaliza"

# --- New row 1111: Escarola ---
$ws.Range("A1111").Value = 11
$ws.Range("B1111").Value = "Vega Monumental Concepción"
$ws.Range("C1111").Value = "Bíobío"
$ws.Range("D1111").Value = 45106
$ws.Range("E1111").Value = 8
$ws.Range("F1111").Value = 100112033
$ws.Range("G1111").Value = "Lechuga"
$ws.Range("H1111").Value = "Escarola"
$ws.Range("I1111").Value = "Primera"
$ws.Range("J1111").Value = 100
$ws.Range("K1111").Value = 6500
$ws.Range("L1111").Value = 7000
$ws.Range("M1111").Value = 6750
$ws.Range("N1111").Value = "$/caja 15 unidades"
$ws.Range("O1111").Value = "Región de Coquimbo"
$ws.Range("P1111").Value = 450
$ws.Range("Q1111").Value = 15
$ws.Range("R1111").Value = "Hortaliza"

# --- New row 1112: Marina ---
$ws.Range("A1112").Value = 11
$ws.Range("B1112").Value = "Vega Monumental Concepción"
$ws.Range("C1112").Value = "Bíobío"
$ws.Range("D1112").Value = 45106
$ws.Range("E1112").Value = 8
$ws.Range("F1112").Value = 100112033
$ws.Range("G1112").Value = "Lechuga"
$ws.Range("H1112").Value = "Marina"
$ws.Range("I1112").Value = "Primera"
$ws.Range("J1112").Value = 100
$ws.Range("K1112").Value = 6000
$ws.Range("L1112").Value = 6500
$ws.Range("M1112").Value = 6250
$ws.Range("N1112").Value = "$/caja 15 unidades"
$ws.Range("O1112").Value = "Región Metropolitana"
$ws.Range("P1112").Value = 417
$ws.Range("Q1112").Value = 15
$ws.Range("R1112").Value = "Hortaliza"
